$wb = $excel.ActiveWorkbook

# --- "About" sheet: bump the "Source:" date to the new update date ---
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Range("C1").Value = 45379

# --- "FPIEBP" sheet: update "hard coal" balancing priorities ---
# columns: B = production, C = imports, D = exports
$wsData = $wb.Worksheets.Item("FPIEBP")
$wsData.Range("B3").Value = 1
$wsData.Range("C3").Value = 3
$wsData.Range("D3").Value = 2

# --- restore view/selection state as last touched by the author ---
$wsAbout.Activate()
$wsAbout.Range("A6").Select()

$wsData.Activate()
$wsData.Range("E3").Select()

$wb.Save()
